$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove row 13 (duplicate sample row tied to the removed local-path shared string).
$ws.Rows.Item(13).Delete()

# 2) Update computed feature values in columns G, H, M, N for rows 2-12
#    (bwArea / entropyVal / sobelArea / cannyArea recomputed for the evaluation stage).
$ws.Range("G2").Value = 201.40649230585709
$ws.Range("H2").Value = 7428.8639984550027
$ws.Range("M2").Value = 25.137841923673367
$ws.Range("N2").Value = 31.7741978959249

$ws.Range("G3").Value = 526.32095336914063
$ws.Range("H3").Value = 1290.0200164768814
$ws.Range("M3").Value = 0.092998504638671875
$ws.Range("N3").Value = 0.10787391662597655

$ws.Range("G4").Value = 626.29904275412082
$ws.Range("H4").Value = 1583.989477721698
$ws.Range("M4").Value = 0.090907815141589188
$ws.Range("N4").Value = 0.10319945979501269

$ws.Range("G5").Value = 904.75320649092976
$ws.Range("H5").Value = 4364.4672557587101
$ws.Range("M5").Value = 2.3932681405895697
$ws.Range("N5").Value = 2.6863555839002267

$ws.Range("G6").Value = 508.76069958847734
$ws.Range("H6").Value = 3506.4597973925534
$ws.Range("M6").Value = 2.6215157750342932
$ws.Range("N6").Value = 2.8948336762688616

$ws.Range("G7").Value = 883.41454399956604
$ws.Range("H7").Value = 5731.9277514401228
$ws.Range("M7").Value = 4.0737847222222223
$ws.Range("N7").Value = 4.1797756618923616

$ws.Range("G8").Value = 839.32014371141975
$ws.Range("H8").Value = 5197.8138623771529
$ws.Range("M8").Value = 3.6155840084876538
$ws.Range("N8").Value = 3.9989149305555554

$ws.Range("G9").Value = 755.87836371527771
$ws.Range("H9").Value = 4490.9417711155975
$ws.Range("M9").Value = 1.161084587191358
$ws.Range("N9").Value = 1.2592110339506171

$ws.Range("G10").Value = 869.42730034722229
$ws.Range("H10").Value = 4692.8090819823165
$ws.Range("M10").Value = 1.1994466145833333
$ws.Range("N10").Value = 1.4099772135416668

$ws.Range("G11").Value = 270.60474537037038
$ws.Range("H11").Value = 6580.9302308018796
$ws.Range("M11").Value = 7.590850453317902
$ws.Range("N11").Value = 8.6450737847222232

$ws.Range("G12").Value = 682.20310443402434
$ws.Range("H12").Value = 5285.9895652530304
$ws.Range("M12").Value = 6.8246005059230015
$ws.Range("N12").Value = 7.4371195294504773

# 3) Resize columns A, M, N to fit the regenerated data.
$ws.Columns.Item(1).ColumnWidth = 58.166666666666664
$ws.Columns.Item(13).ColumnWidth = 12.833333333333332
$ws.Columns.Item(14).ColumnWidth = 11.833333333333332
